$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ProjectSchedule")

# --- Update phase/task progress percentages ---
$ws.Range("D10").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("D16").Value = 1
$ws.Range("D21").Value = 0.8
$ws.Range("D22").Value = 1
$ws.Range("D23").Value = 0.8
$ws.Range("D24").Value = 1
$ws.Range("D25").Value = 0.8
$ws.Range("D28").Value = 0.2

# --- Rename task / reassign owner for the new software process model step ---
$ws.Range("B20").Value = "Sprint 4 - Step 4"

# --- Row 27: renamed task ---
$ws.Range("B27").Value = "Test Report"

# --- Row 30: new task added to the schedule ---
$ws.Range("B30").Value = "Launch to Heroku "
$ws.Range("C30").Value = "Jaisal Friedman"

# --- Row 23: replaced task + reassigned owner ---
$ws.Range("B23").Value = "Request to Connect with Mentor Dashboard"
$ws.Range("C23").Value = "Daniel Waston"

# --- Update the visible scroll/selection state to match the latest edits ---
$ws.Application.ActiveWindow.ScrollRow = 21
$ws.Range("D29").Select()
